$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.477.05"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.57%  "
$ws.Range("D3").Value = "'1.839.85"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.94%  "
$ws.Range("D5").Value = "'243.30"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.41%  "
$ws.Range("D6").Value = "'0.6274"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.14%  "
$ws.Range("B8").Value = "Dogecoin"
$ws.Range("C8").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D8").Value = "'0.07416"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.61%  "
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D9").Value = "'0.2939"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.55%  "
$ws.Range("D10").Value = "'23.39"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.34%  "
$ws.Range("D11").Value = "'0.07643"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.20%  "
$ws.Range("D12").Value = "'1.829.68"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'5.017"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.15%  "
$ws.Range("D14").Value = "'0.6766"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.07%  "
$ws.Range("D15").Value = "'83.34"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.19%  "
$ws.Range("D16").Value = "'0.000009398"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.23%  "
$ws.Range("D17").Value = "'5.898"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.15%  "
$ws.Range("D18").Value = "'29.446.75"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Value = "'2.083.28"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").Value = "'238.03"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.48%  "
$ws.Range("D21").Value = "'12.54"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.03%  "
$ws.Range("D22").Value = "'1.000"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.29%  "
$ws.Range("D23").Value = "'7.350"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +2.12%  "
$ws.Range("D24").Value = "'0.9995"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.12%  "
$ws.Range("D25").Value = "'158.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.89%  "
$ws.Range("D26").Value = "'0.1415"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.58%  "
$ws.Range("D27").Value = "'8.500"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.49%  "
$ws.Range("D28").Value = "'17.77"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("D29").Value = "'0.06091"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +8.31%  "
$ws.Range("D30").Value = "'1.496"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("D31").Value = "'1.231"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.77%  "
$ws.Range("D32").Value = "'4.095"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.76%  "
$ws.Range("D33").Value = "'4.112"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "
$ws.Range("D34").Value = "'1.869"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.74%  "
$ws.Range("D35").Value = "'1.144"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.07%  "
$ws.Range("D36").Value = "'0.7256"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.62%  "
$ws.Range("D37").Value = "'2.613"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.87%  "
$ws.Range("D38").Value = "'2.877"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.69%  "
$ws.Range("D39").Value = "'1.219.69"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.13%  "
$ws.Range("E40").Value = "  -1.47%  "
$ws.Range("D41").Value = "'6.309"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.64%  "
$ws.Range("D42").Value = "'0.9132"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.06%  "
$ws.Range("E43").Value = "  -1.15%  "
$ws.Range("D44").Value = "'1.997.38"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'101.97"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "'65.64"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "'0.5065"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.31%  "
$ws.Range("D48").Value = "'0.00000000120"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.10%  "
$ws.Range("D49").Value = "'9.260"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +1.07%  "
$ws.Range("D50").Value = "'0.4058"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.07%  "
$ws.Range("D51").Value = "'0.1141"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.51%  "
